$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.219.40"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "1.796.99"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'314.29"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  +1.67%  "
$ws.Range("D8").Value = "'0.3817"
$ws.Range("E8").Value = "  -3.06%  "
$ws.Range("D9").Value = "'0.07921"
$ws.Range("E9").Value = "  -3.39%  "
$ws.Range("D10").Value = "'41.67"
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("D11").Value = "'1.097"
$ws.Range("E11").Value = "  -1.57%  "
$ws.Range("D12").Value = "'6.271"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").Value = "'1.001"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "'20.55"
$ws.Range("E14").Value = "  -2.76%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'7.264"
$ws.Range("E15").Value = "  -3.97%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.788.46"
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("D17").Value = "'93.07"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "'0.00001083"
$ws.Range("E18").Value = "  -3.53%  "
$ws.Range("D19").Value = "'0.06546"
$ws.Range("E19").Value = "  -1.78%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "'17.30"
$ws.Range("E21").Value = "  -3.04%  "
$ws.Range("D22").Value = "'5.951"
$ws.Range("E22").Value = "  -2.66%  "
$ws.Range("D23").Value = "28.274.39"
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("D24").Value = "'11.13"
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").Value = "'160.51"
$ws.Range("E26").Value = "  +2.12%  "
$ws.Range("D27").Value = "'20.46"
$ws.Range("E27").Value = "  -4.03%  "
$ws.Range("D28").Value = "2.000.72"
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("D29").Value = "'2.337"
$ws.Range("E29").Value = "  -3.29%  "
$ws.Range("D30").Value = "'123.23"
$ws.Range("E30").Value = "  -2.60%  "
$ws.Range("E31").Value = "  -1.75%  "
$ws.Range("E32").Value = "  -5.37%  "
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").Value = "'5.559"
$ws.Range("E34").Value = "  -3.66%  "
$ws.Range("D35").Value = "'0.07307"
$ws.Range("E35").Value = "  +3.57%  "
$ws.Range("D36").Value = "'12.20"
$ws.Range("E36").Value = "  +8.08%  "
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("D38").Value = "'0.2145"
$ws.Range("E38").Value = "  -3.96%  "
$ws.Range("D39").Value = "'5.076"
$ws.Range("E39").Value = "  -3.37%  "
$ws.Range("D40").Value = "'8.612"
$ws.Range("E40").Value = "  -2.08%  "
$ws.Range("D41").Value = "'0.6169"
$ws.Range("E41").Value = "  -3.20%  "
$ws.Range("D42").Value = "'1.161"
$ws.Range("E42").Value = "  -1.84%  "
$ws.Range("D43").Value = "'1.370"
$ws.Range("E43").Value = "  -2.19%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'13.23"
$ws.Range("E44").Value = "  -3.17%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.6029"
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("D46").Value = "'3.783"
$ws.Range("E46").Value = "  +1.16%  "
$ws.Range("D47").Value = "'127.36"
$ws.Range("E47").Value = "  +1.66%  "
$ws.Range("D48").Value = "'1.231"
$ws.Range("E48").Value = "  +2.90%  "
$ws.Range("D49").Value = "'1.922"
$ws.Range("E49").Value = "  -3.55%  "
$ws.Range("D50").Value = "'0.06778"
$ws.Range("E50").Value = "  -2.39%  "
$ws.Range("D51").Value = "'73.10"
$ws.Range("E51").Value = "  -1.65%  "

Write-Host "Applied cryptos update"
